$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Remove the two oldest fixtures (West Ham United match + MU Women v
#        Tottenham Hotspur Women), which have already been played and are
#        dropped from the list. Locate them by content so the edit is not
#        dependent on a fixed row offset.
$westHam = $ws.Range("A1:A30").Find("West Ham United")
if ($westHam -ne $null) {
    $westHam.EntireRow.Delete() | Out-Null
}
$tottenhamWomen = $ws.Range("A1:A30").Find("Tottenham Hotspur Women")
if ($tottenhamWomen -ne $null) {
    $tottenhamWomen.EntireRow.Delete() | Out-Null
}

# --- 2) Insert the newly-announced fixture (Brighton and Hove Albion, 15 Feb)
#        directly below the "Southampton" row.
$southampton = $ws.Range("A1:A30").Find("Southampton")
$newRow = $southampton.Row + 1
$ws.Rows.Item($newRow).Insert()
$ws.Cells.Item($newRow, 1).Value = "Manchester United v Brighton and Hove Albion "
$ws.Cells.Item($newRow, 2).Value = "15 FebTue20:15"

# --- 3) The Tottenham Hotspur home fixture kicked off time moved from
#        15:00 to 17:30 (TV rescheduling).
$tottenham = $ws.Range("A1:A30").Find("Tottenham Hotspur  ")
$ws.Cells.Item($tottenham.Row, 2).Value = "12 MarSat17:30"
